$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the greeting cell from "Good Morning" to "GIT UPDATE"
$ws.Range("E8").Value = "GIT UPDATE"

# Select the edited cell so the sheetView selection matches the diff
$ws.Range("E8").Select()
